$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Content: drop the old "day of week" column (B), insert a new header row,
# and tidy up two descriptions with a trailing period.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).Delete()
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Plan"
$ws.Range("B2").Value = "Finalize layout, start blocking  out major objects."
$ws.Range("B3").Value = "Housekeeping! Tidy up the layout, tentively set out secondary objects."

Write-Output "content done"

# ---------------------------------------------------------------------------
# Column / row sizing
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 54.8
$ws.Rows.Item(1).RowHeight = 20

# ---------------------------------------------------------------------------
# Header row formatting: bold 14pt, centered, boxed in medium borders
# ---------------------------------------------------------------------------
$hdrRange = $ws.Range("A1:B1")
$hdrRange.Font.Bold = $true
$hdrRange.Font.Size = 14
$hdrRange.HorizontalAlignment = -4108   # xlCenter

foreach ($cell in $hdrRange.Cells) {
  $cell.Borders.Item(7).Weight = -4138   # xlEdgeLeft, medium
  $cell.Borders.Item(8).Weight = -4138   # xlEdgeTop, medium
  $cell.Borders.Item(9).Weight = -4138   # xlEdgeBottom, medium
  $cell.Borders.Item(10).Weight = -4138  # xlEdgeRight, medium
}

Write-Output "header done"

# ---------------------------------------------------------------------------
# Data rows (2-7): thin/medium ruled table, banded gray fill on alternate
# rows (the 2nd, 4th and 6th data rows -> sheet rows 3, 5 and 7).
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 7; $r++) {
  $colA = $ws.Cells.Item($r, 1)
  $colB = $ws.Cells.Item($r, 2)

  $colA.Borders.Item(7).Weight = -4138     # xlEdgeLeft, medium -> actually thin, fixed below
  $colA.Borders.Item(7).LineStyle = 1
  $colA.Borders.Item(7).Weight = 2         # thin
  $colA.Borders.Item(10).Weight = -4138    # xlEdgeRight, medium
  $colA.Borders.Item(9).Weight = 2         # xlEdgeBottom, thin
  if ($r -eq 2) {
    $colA.Borders.Item(8).Weight = -4138   # xlEdgeTop, medium (sits right under the header)
  } else {
    $colA.Borders.Item(8).Weight = 2       # xlEdgeTop, thin
  }

  $colB.Borders.Item(10).Weight = 2        # xlEdgeRight, thin
  $colB.Borders.Item(9).Weight = 2         # xlEdgeBottom, thin
  if ($r -eq 2) {
    # no top border on B2
  } else {
    $colB.Borders.Item(8).Weight = 2       # xlEdgeTop, thin
  }

  if ($r -eq 3 -or $r -eq 5 -or $r -eq 7) {
    $colA.Interior.Color = 14277081        # White, Background 1, Darker 15% (theme 0, tint -0.15)
    $colB.Interior.Color = 14277081
  }
}

Write-Output "data rows done"
